$wb = $excel.ActiveWorkbook

# --- Sheet: time_variants ---
$ws2 = $wb.Worksheets.Item("time_variants")

# Delete row 8 (int_perc_improve_dst) entirely; rows below shift up by one
$ws2.Rows.Item(8).Delete() | Out-Null

# Delete columns X:Z (scenario_7, scenario_8, scenario_9)
$ws2.Range("X1:Z1").EntireColumn.Delete() | Out-Null

Write-Host "after deletes"
Write-Host $ws2.Range("A8").Text
Write-Host $ws2.Range("A16").Text
Write-Host $ws2.UsedRange.Address()
